# [yerim] Ap 와 Actor data 분리, monsterskillbase 생성
#
# actor_rsc sheet: populate the (previously empty) C column — "actor_rsc_prefab" —
# for the two existing rows, splitting the player/enemy actor resource prefabs
# apart instead of leaving them blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("actor_rsc")

$ws.Range("C4").Value = "ActorPlayer"
$ws.Range("C5").Value = "ActorEnemy"

# Leave the selection where the author ended up after filling these cells in.
$ws.Range("C5").Select()
